$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.896.07'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.443.51'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.88%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '583.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.34%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.70'
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.477'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.92%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.63'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.23%  '

$ws.Range("E10").Value = '  +2.39%  '

$ws.Range("E11").Value = '  +2.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.036.71'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '29.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.99%  '

$ws.Range("E14").Value = '  -0.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.445.57'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.08%  '

$ws.Range("E16").Value = '  +2.09%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '62.888.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.94%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.23'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.32%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.49%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.30'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.85%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '394.38'
$ws.Range("D21").Style = "Normal"

$ws.Range("B22").Value = 'Litecoin'
$ws.Range("C22").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.16%  '

$ws.Range("B23").Value = 'Polygon'
$ws.Range("C23").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.562'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.41%  '

$ws.Range("E24").Value = '  +0.13%  '

$ws.Range("E25").Value = '  +4.61%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.586.49'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.99%  '

$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.35%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.998'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.18'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.05%  '

$ws.Range("E31").Value = '  +6.61%  '

$ws.Range("E32").Value = '  +1.30%  '

$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.74'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.37%  '

$ws.Range("E35").Value = '  +7.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.07'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.39%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.59'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '168.06'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '30.61'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +18.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.477.00'
$ws.Range("D40").Style = "Normal"

$ws.Range("E41").Value = '  +0.96%  '

$ws.Range("E42").Value = '  +1.57%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.87'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.27%  '

$ws.Range("E44").Value = '  +3.15%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.73'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.61%  '

$ws.Range("E46").Value = '  +8.42%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.518.21'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.63%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.60'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.51%  '

$ws.Range("E49").Value = '  +1.94%  '

$ws.Range("E50").Value = '  +0.04%  '

$ws.Range("E51").Value = '  +4.54%  '
